# Insert a new data row at row 507, shifting existing rows 507:607 down to 508:608.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(507).Insert()

$ws.Range("A507").Value = 10
$ws.Range("B507").Value = "Vega Modelo de Temuco"
$ws.Range("C507").Value = "La Araucanía"
$ws.Range("D507").Value = 45275
$ws.Range("E507").Value = 9
$ws.Range("F507").Value = 100112009
$ws.Range("G507").Value = "Acelga"
$ws.Range("H507").Value = "Sin especificar"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 150
$ws.Range("K507").Value = 10000
$ws.Range("L507").Value = 10000
$ws.Range("M507").Value = 10000
$ws.Range("N507").Value = "$/docena de atados (12 kilos)"
$ws.Range("O507").Value = "Región de La Araucanía"
$ws.Range("P507").Value = 833
$ws.Range("Q507").Value = 12
$ws.Range("R507").Value = "Hortaliza"
